$wb = $excel.ActiveWorkbook

# --- Per-technology cost sheets (Operating / Maintenance / Capital / Total) ---
# The "Electric_boiler" technology row (row 1) is removed from each of these
# four sheets. Deleting the row shifts the remaining technology rows up by
# one and drops the now-unused "Electric_boiler" shared string, which is
# exactly the shared-string / dimension change described in the diff.

$techSheetValues = @{
    "Operating_cost_per_technology"    = @(72799.921433088602, 84453.791263466919, 0, 0, 0, 0)
    "Maintenance_cost_per_technology"  = @(13912.064985863608, 8820.7293097398433, 0, 0, 14160.972742735672, 0)
    "Capital_cost_per_technology"      = @(14531.307619497124, 32044.611431009191, 0, 0, 34395.988641924261, 0)
    "Total_cost_per_technology"        = @(101243.29403844933, 125319.13200421595, 0, 0, 48556.961384659931, 0)
}

foreach ($sheetName in $techSheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Rows.Item(1).Delete()

    $values = $techSheetValues[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($i + 1, 2).Value = $values[$i]
    }
}

# --- Grid cost sheets ---
$wb.Worksheets.Item("Operating_cost_grid").Range("A1").Value = 195246.80175913486
$wb.Worksheets.Item("Total_cost_grid").Range("A1").Value = 195246.80175913486

# --- Storage cost sheets (Elec / Heat rows; labels renumber automatically) ---
$wb.Worksheets.Item("Capital_cost_per_storage").Range("B2").Value = 3829.1911649368767
$wb.Worksheets.Item("Total_cost_per_storage").Range("B2").Value = 3829.1911649368767

# --- Income via exports ---
$wb.Worksheets.Item("Income_via_exports").Range("A1").Value = 12612.424068058446
